# Fruta / hortaliza, semanal
#
# Insert two new weekly price records (rows 191-192) into the daily logic
# sheet for "Fruta, Terminal La Palmera de La Serena - Ciruela", pushing the
# existing rows 191-229 down to 193-231.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 191, shifting everything below down by two rows.
$ws.Rows("191:192").Insert()

# Common columns (A-J, Q, T) shared by every "Angeleno" bins record in this
# block.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100103
$producto  = "Frutos de hueso (carozo)"
$categoriaId = 100103002
$categoria = "Ciruela"
$variedad  = "Angeleno"
$unidad    = "$/bins (450 kilos)"
$origen    = "Región de O'Higgins"
$kgUnidad  = 450

function Set-CiruelaRow {
    param($Row, $Fecha, $Calidad, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $PrecioKg)

    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $tipo
    $ws.Cells.Item($Row, 7).Value = $productoId
    $ws.Cells.Item($Row, 8).Value = $producto
    $ws.Cells.Item($Row, 9).Value = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMinimo
    $ws.Cells.Item($Row, 15).Value = $PrecioMaximo
    $ws.Cells.Item($Row, 16).Value = $PrecioPromedio
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad
}

# New row 191: Angeleno / Especial
Set-CiruelaRow 191 44641 "Especial" 10 230000 240000 235000 522

# New row 192: Angeleno / Primera
Set-CiruelaRow 192 44641 "Primera" 10 210000 220000 215000 478
